# Updates the crypto price/volume table (columns D "Price" and E "Volume(1h)")
# on Sheet1, rows 2-51, to the latest scraped values from the GitHub Actions
# refresh run. Row 1 is the header and is left untouched.
#
# Column D is stored as literal text in the source data (some prices use
# "." as a thousands separator, e.g. "58.022.63", which is not a valid
# number) so every write below is kept as text. For the subset of new D
# values that *are* plain decimals (e.g. "518.14"), Excel would otherwise
# silently reinterpret the text as a Number on assignment; forcing the
# cell's NumberFormat to "@" (Text) first preserves the exact literal
# string, matching how the rest of the column is already stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '58.022.63'; E = '  -1.50%  ' },
    @{ Row = 3; D = '2.469.24'; E = '  -2.27%  ' },
    @{ Row = 4; E = '  -0.20%  ' },
    @{ Row = 5; D = '518.14'; E = '  -3.52%  ' },
    @{ Row = 6; D = '131.43'; E = '  -4.13%  ' },
    @{ Row = 7; E = '  +0.02%  ' },
    @{ Row = 8; D = '0.558'; E = '  -1.88%  ' },
    @{ Row = 9; D = '0.0992'; E = '  -2.20%  ' },
    @{ Row = 10; E = '  -0.42%  ' },
    @{ Row = 11; E = '  +0.06%  ' },
    @{ Row = 12; E = '  -1.64%  ' },
    @{ Row = 13; D = '2.906.40'; E = '  -2.10%  ' },
    @{ Row = 14; D = '57.958.71'; E = '  -1.65%  ' },
    @{ Row = 15; D = '22.30'; E = '  -3.31%  ' },
    @{ Row = 16; D = '0.0000137'; E = '  -2.12%  ' },
    @{ Row = 17; D = '2.467.83'; E = '  -1.85%  ' },
    @{ Row = 18; E = '  -3.49%  ' },
    @{ Row = 19; E = '  -2.37%  ' },
    @{ Row = 20; D = '320.38'; E = '  -1.05%  ' },
    @{ Row = 21; E = '  +0.12%  ' },
    @{ Row = 22; D = '5.75'; E = '  -3.50%  ' },
    @{ Row = 23; D = '64.06'; E = '  -2.54%  ' },
    @{ Row = 24; D = '0.412'; E = '  -2.52%  ' },
    @{ Row = 25; D = '0.999'; E = '  +0.22%  ' },
    @{ Row = 26; E = '  -2.97%  ' },
    @{ Row = 27; D = '7.34'; E = '  -2.64%  ' },
    @{ Row = 28; D = '0.0₃0750'; E = '  -3.09%  ' },
    @{ Row = 29; D = '166.10'; E = '  -0.97%  ' },
    @{ Row = 30; E = '  -4.41%  ' },
    @{ Row = 31; D = '6.30'; E = '  -6.24%  ' },
    @{ Row = 32; E = '  -2.03%  ' },
    @{ Row = 33; D = '0.998'; E = '  +0.12%  ' },
    @{ Row = 34; D = '0.999'; E = '  +0.16%  ' },
    @{ Row = 35; D = '18.02'; E = '  -2.25%  ' },
    @{ Row = 36; E = '  -10.55%  ' },
    @{ Row = 37; E = '  -3.41%  ' },
    @{ Row = 38; E = '  -4.63%  ' },
    @{ Row = 39; D = '0.789'; E = '  -3.36%  ' },
    @{ Row = 40; D = '3.46'; E = '  -4.67%  ' },
    @{ Row = 41; D = '273.10'; E = '  -3.97%  ' },
    @{ Row = 42; D = '4.99'; E = '  -2.87%  ' },
    @{ Row = 43; D = '0.592'; E = '  -2.49%  ' },
    @{ Row = 44; D = '126.56'; E = '  -3.91%  ' },
    @{ Row = 45; D = '0.0905'; E = '  -2.16%  ' },
    @{ Row = 46; D = '0.0488'; E = '  -4.11%  ' },
    @{ Row = 47; E = '  -3.25%  ' },
    @{ Row = 48; E = '  -2.01%  ' },
    @{ Row = 49; D = '1.731.77'; E = '  -1.96%  ' },
    @{ Row = 50; E = '  -1.44%  ' },
    @{ Row = 51; E = '  -1.10%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D" + $u.Row)
        # Plain-decimal strings (e.g. "518.14") would otherwise be silently
        # reinterpreted by Excel as a Number on assignment; values using "."
        # as a thousands separator (e.g. "58.022.63") already fail that
        # auto-detection on their own, so only guard the ones that match.
        if ($u.D -match '^-?\d+(\.\d+)?$') {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
